# Auto update Excel log
# Appends newly detected sensor events to the ALERTS and mmWave logs.

$wb = $excel.ActiveWorkbook

function Add-LogRow {
    param($ws, $row, $date, $timestamp, $hour, $location, $value, $status)

    # Column A holds a date-looking string (e.g. "2026-02-01"). Force the
    # cell to Text format before assigning it so Excel stores the literal
    # string instead of auto-converting it to a date serial number (as the
    # rest of the log does), then reset the style back to Normal so no
    # extra number formatting is left applied to the cell.
    $cellA = $ws.Range("A" + $row)
    $cellA.NumberFormat = "@"
    $cellA.Value = $date
    $cellA.Style = "Normal"

    $ws.Range("B" + $row).Value = $timestamp
    $ws.Range("C" + $row).Value = $hour
    $ws.Range("D" + $row).Value = $location
    $ws.Range("E" + $row).Value = $value
    $ws.Range("F" + $row).Value = $status
}

# ---- ALERTS sheet: two new FALL_DETECTED critical alerts ----
$alerts = $wb.Worksheets.Item("ALERTS")
Add-LogRow $alerts 8 "2026-02-01" "00:21:33" "00:00" "Living Room" "CRITICAL" "FALL_DETECTED"
Add-LogRow $alerts 9 "2026-02-01" "00:21:37" "00:00" "Living Room" "CRITICAL" "FALL_DETECTED"

# ---- mmWave sheet: three new presence/motion events ----
$mmwave = $wb.Worksheets.Item("mmWave")
Add-LogRow $mmwave 53 "2026-02-01" "00:21:52" "00:00" "Living Room" "NO_MOTION_DETECTED" "Inactive"
Add-LogRow $mmwave 54 "2026-02-01" "00:22:02" "00:00" "Living Room" "PRESENCE_DETECTED"  "Active"
Add-LogRow $mmwave 55 "2026-02-01" "00:22:20" "00:00" "Living Room" "PRESENCE_DETECTED"  "Active"
